# English translation of sheet names
$wb = $excel.ActiveWorkbook

# Rename sheets to their English equivalents (order matches existing tab order)
$wb.Worksheets.Item(1).Name = "Samples1_impl"
$wb.Worksheets.Item(2).Name = "Samples2_impl"
$wb.Worksheets.Item(3).Name = "Models_impl"
$wb.Worksheets.Item(4).Name = "AutoTruncations_impl"
$wb.Worksheets.Item(5).Name = "SystematicTruncations_impl"
$wb.Worksheets.Item(6).Name = "SpecialTroncations_expl"

# Move the active/selected tab from the 1st sheet to the 5th sheet
$wb.Worksheets.Item(5).Activate()
